$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header "rotacao_poste" is inserted; the former K column
# (header "lista_abaco" plus its bracketed coordinate-list values) moves to L.
$ws.Range("L1").Value = "lista_abaco"
$ws.Range("K1").Value = "rotacao_poste"

$ws.Range("L2").Value = "[(-1.0, 0.0), (-1.0, 217.143), (160.0, 217.143), (160.0, 0.0), (-1.0, 0.0)]]"
$ws.Range("K2").Value = "topo1"

$ws.Range("K3").Value = "bissetriz2"
$ws.Range("L3").Value = "[(-1.0, 0.0), (-1.0, 217.143), (160.0, 217.143), (160.0, 0.0), (-1.0, 0.0)]]"

$ws.Range("K6").Value = "bissetriz2"

# Update the selected cell to match the saved view state.
$ws.Range("K3").Select()
